$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C; this shifts the existing
# column C (the "Jun_10" column) two places to the right, to column E.
$ws.Range("C:D").Insert()

# The old B1 header ("Jun_13") needs to move out to the new D1 slot,
# and B1/C1 become brand-new headers ("Jun_17"/"Jun_15").
$ws.Range("D1").Value = $ws.Range("B1").Value()
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the two new columns (C & D) with the same "UN" marker used
# throughout column B, for every data row (2-27).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
